$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Sheet1: Forecast Comparison
$ws1.Range("L2").Value = 1.02
$ws1.Range("L3").Value = 1.03

$ws1.Range("H4").Value = 19.55
$ws1.Range("L4").Value = 1.18

$ws1.Range("H5").Value = 18.55
$ws1.Range("L5").Value = 0.84

$ws1.Range("H6").Value = 17.55
$ws1.Range("L6").Value = 1.15

$ws1.Range("H7").Value = 16.55
$ws1.Range("L7").Value = 1.04

$ws1.Range("H8").Value = 17.1
$ws1.Range("L8").Value = 1.11

$ws1.Range("H9").Value = 16.1
$ws1.Range("L9").Value = 1.12

$ws1.Range("H10").Value = 13.73
$ws1.Range("L10").Value = 0.8

$ws1.Range("H11").Value = 11.67
$ws1.Range("L11").Value = 0.9

$ws1.Range("H12").Value = 10.67
$ws1.Range("L12").Value = 0.93

$ws1.Range("H13").Value = 9.67
$ws1.Range("L13").Value = 0.89

$ws1.Range("H14").Value = 8
$ws1.Range("L14").Value = 0.97

$ws1.Range("H15").Value = 7
$ws1.Range("L15").Value = 0.88

$ws1.Range("H16").Value = 6.5
$ws1.Range("L16").Value = 1.08

$ws1.Range("H17").Value = 6.6
$ws1.Range("L17").Value = 0.93

# Sheet2: Summary
# Values here are stored as text (not numbers) in the source file, so force
# the cell format to Text before assigning, then restore the original
# (unstyled / "Normal") cell style so no stray formatting is introduced.
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "36"
$ws2.Range("B9").Style = "Normal"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "17"
$ws2.Range("B10").Style = "Normal"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "3"
$ws2.Range("B12").Style = "Normal"
